$wb = $excel.ActiveWorkbook

# --- Sheet: y_fitted_on_begin_2016 ---
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$values1 = @(
    -0.2610738699168351,
    15.82610396925855,
    15.49677399888529,
    15.18464152532538,
    14.82905189465676,
    13.33143331283013,
    13.23543880064474,
    13.08456824468752,
    12.94530615472575,
    12.73590430331101,
    12.82120806736439,
    12.49019561909548,
    11.85739694843217,
    11.38072621320181,
    11.83377893261917,
    11.72694333012629,
    11.64341911215989,
    11.59262719038455,
    11.41698376269012,
    11.4893894322177
)

$row = 2
foreach ($v in $values1) {
    $ws1.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Sheet: y_fitted_on_begin_2021 ---
$ws2 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$values2 = @(
    -0.2488829571298799,
    15.83829462197516,
    15.51793299211485,
    15.20317715015738,
    14.85570323437172,
    13.55894593152401,
    13.19136693218302,
    13.08252766666925,
    12.9362097408316,
    12.74004512393225,
    12.77168362779146,
    12.52064423597761,
    11.93166535821961,
    11.42218207336993,
    11.71481225924756,
    11.72603591405079,
    11.62433608847502,
    11.56992258222618,
    11.41677786528379,
    11.4893894322177,
    10.9486794529308,
    11.04244184491913,
    11.12900142796406,
    10.81593349450332,
    10.98316395120389
)

$row = 2
foreach ($v in $values2) {
    $ws2.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}
